# Schulferien - unnötige Zeilen gelöscht
# The "Hinweise"/footnote block in rows 20-28 (columns A:G) contained a bunch of
# explanatory text (Stand-Datum, Fußnoten, Erläuterungen) that is no longer
# needed. Remove its contents while keeping the existing cell formatting
# (borders/fonts/merges) intact, and leave the selection on the cleared range,
# matching the author's save state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the values/text of the now-unnecessary explanatory rows, keeping styles.
$ws.Range("A20:G28").ClearContents()

# Reflect the selection left behind after deleting the content.
$ws.Range("A18:G28").Select()
